$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.3030883333333333
$ws.Range("H2").Value = 0.909265
$ws.Range("I2").Value = 0.5850568929085261
$ws.Range("J2").Value = 0.585056892908526
$ws.Range("M2").Value = 44.68160133333333
$ws.Range("N2").Value = 134.044804
$ws.Range("O2").Value = 0.9072345081554035
$ws.Range("P2").Value = 0.9072345081554034
$ws.Range("Q2").Value = 13.54247207878445
$ws.Range("R2").Value = 121.88224870906
$ws.Range("S2").Value = 0.5307838024807953
$ws.Range("T2").Value = 0.5307838024807952

# Row 3
$ws.Range("G3").Value = 0.3030883333333333
$ws.Range("H3").Value = 0.909265
$ws.Range("I3").Value = 0.5850568929085261
$ws.Range("J3").Value = 0.585056892908526
$ws.Range("O3").Value = 0.08277108896415035
$ws.Range("P3").Value = 0.08277108896415035
$ws.Range("Q3").Value = 1.235540702157222
$ws.Range("R3").Value = 11.119866319415
$ws.Range("S3").Value = 0.048425796132021
$ws.Range("T3").Value = 0.04842579613202099

# Row 4
$ws.Range("G4").Value = 0.3030883333333333
$ws.Range("H4").Value = 0.909265
$ws.Range("I4").Value = 0.5850568929085261
$ws.Range("J4").Value = 0.585056892908526
$ws.Range("M4").Value = 0.3223466666666667
$ws.Range("N4").Value = 0.96704
$ws.Range("O4").Value = 0.006545065773430512
$ws.Range("P4").Value = 0.006545065773430512
$ws.Range("Q4").Value = 0.09769951395555557
$ws.Range("R4").Value = 0.8792956256
$ws.Range("S4").Value = 0.003829235845285195
$ws.Range("T4").Value = 0.003829235845285194

# Row 5
$ws.Range("G5").Value = 0.3030883333333333
$ws.Range("H5").Value = 0.909265
$ws.Range("I5").Value = 0.5850568929085261
$ws.Range("J5").Value = 0.585056892908526
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.169881
$ws.Range("N5").Value = 0.5096430000000001
$ws.Range("O5").Value = 0.003449337107015684
$ws.Range("P5").Value = 0.003449337107015683
$ws.Range("Q5").Value = 0.05148894915500001
$ws.Range("R5").Value = 0.4634005423950001
$ws.Range("S5").Value = 0.00201805845042468
$ws.Range("T5").Value = 0.00201805845042468

# Row 6
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.214961
$ws.Range("H6").Value = 0.644883
$ws.Range("I6").Value = 0.4149431070914739
$ws.Range("J6").Value = 0.4149431070914739
$ws.Range("M6").Value = 44.68160133333333
$ws.Range("N6").Value = 134.044804
$ws.Range("O6").Value = 0.9072345081554035
$ws.Range("P6").Value = 0.9072345081554034
$ws.Range("Q6").Value = 9.604801704214665
$ws.Range("R6").Value = 86.443215337932
$ws.Range("S6").Value = 0.3764507056746082
$ws.Range("T6").Value = 0.3764507056746082

# Row 7
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.214961
$ws.Range("H7").Value = 0.644883
$ws.Range("I7").Value = 0.4149431070914739
$ws.Range("J7").Value = 0.4149431070914739
$ws.Range("O7").Value = 0.08277108896415035
$ws.Range("P7").Value = 0.08277108896415035
$ws.Range("Q7").Value = 0.8762893046903333
$ws.Range("R7").Value = 7.886603742213
$ws.Range("S7").Value = 0.03434529283212935
$ws.Range("T7").Value = 0.03434529283212935

# Row 8
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.214961
$ws.Range("H8").Value = 0.644883
$ws.Range("I8").Value = 0.4149431070914739
$ws.Range("J8").Value = 0.4149431070914739
$ws.Range("M8").Value = 0.3223466666666667
$ws.Range("N8").Value = 0.96704
$ws.Range("O8").Value = 0.006545065773430512
$ws.Range("P8").Value = 0.006545065773430512
$ws.Range("Q8").Value = 0.06929196181333333
$ws.Range("R8").Value = 0.62362765632
$ws.Range("S8").Value = 0.002715829928145317
$ws.Range("T8").Value = 0.002715829928145317

# Row 9
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.214961
$ws.Range("H9").Value = 0.644883
$ws.Range("I9").Value = 0.4149431070914739
$ws.Range("J9").Value = 0.4149431070914739
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.169881
$ws.Range("N9").Value = 0.5096430000000001
$ws.Range("O9").Value = 0.003449337107015684
$ws.Range("P9").Value = 0.003449337107015683
$ws.Range("Q9").Value = 0.036517789641
$ws.Range("R9").Value = 0.328660106769
$ws.Range("S9").Value = 0.001431278656591004
$ws.Range("T9").Value = 0.001431278656591003
